# Applies the "changed numvlmon, and added transnorm parameter" commit.
#
# What actually happens (reverse-engineered from the OOXML diff):
#   1. A brand-new parameter row ("transnorm" / "Transmission normalization
#      factor") is inserted as row 46 of the "Model parameters" sheet. This
#      pushes every existing row from 46..117 down by one (47..118), which is
#      exactly what Excel's native Rows.Insert does - content, number format
#      and styles all travel with the row below.
#   2. The active sheet/tab switches from "Transitions" to "Model parameters",
#      and the selection on "Model parameters" moves to B46 (the label cell of
#      the freshly inserted row) while the frozen pane scrolls to A26.

$wb = $excel.ActiveWorkbook

$wsTransitions = $wb.Worksheets.Item("Transitions")
$wsParams = $wb.Worksheets.Item("Model parameters")

# --- 1. Insert the new "transnorm" parameter row ------------------------
# Inserting a whole row shifts rows 46-117 down to 47-118 and carries their
# formatting/styles along, matching the target workbook exactly.
$wsParams.Rows.Item(46).Insert()

$wsParams.Range("A46").Value = "transnorm"
$wsParams.Range("B46").Value = "Transmission normalization factor"
$wsParams.Range("C46").Value = "(0, 'maxmeta')"
$wsParams.Range("D46").Value = "tot"
$wsParams.Range("E46").Value = "constant"
$wsParams.Range("F46").Value = "const"
$wsParams.Range("G46").Value = 0
$wsParams.Range("H46").Value = 0
$wsParams.Range("I46").Value = 0

# --- 2. Update the active tab / selection --------------------------------
# Previously "Transitions" was the selected tab with selection K25; now
# "Model parameters" is selected, with the cursor on the new row's label.
$wsTransitions.Range("K25").Select() | Out-Null

$wsParams.Activate() | Out-Null
$wsParams.Range("B46").Select() | Out-Null

Write-Output "Inserted transnorm parameter row and updated active sheet/selection"
